$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove old rows 3/4 content entirely (clean slate) ---
$ws.Range("A3:G4").Clear() | Out-Null
$ws.Range("F2").Clear() | Out-Null

# --- Row 1 / Row 2 values: write values in an order that keeps surviving
#     shared strings at their original indices, then introduces new ones in
#     the exact order they first appear in the target workbook ---
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "Email"
$ws.Range("E1").Value = "Password"
$ws.Range("A2").Value = "Mahmoud"
$ws.Range("B2").Value = "ElSharkawy"
$ws.Range("C2").Value = "mahmoud.elsharkawy"
$ws.Range("E2").Value = "m1155150745"

$ws.Range("F1").Value = "MobileNumber"
$ws.Range("G1").Value = "Country"
$ws.Range("G2").Value = "Egypt"
$ws.Range("J2").Value = "adress1"
$ws.Range("K2").Value = "adress2"
$ws.Range("I1").Value = "City"
$ws.Range("I2").Value = "Cairo"
$ws.Range("D1").Value = "EmailDomain"
$ws.Range("D2").Value = "test.com"
$ws.Range("L1").Value = "Status"
$ws.Range("L2").Value = "Enabled"
$ws.Range("H1").Value = "NewCity"
$ws.Range("H2").Value = "Alex"
$ws.Range("J1").Value = "Address1"
$ws.Range("K1").Value = "Address2"

$ws.Range("F2").Value = 12345678901

# --- Styling ---
# Header row (A1:L1) -> style index 3 (fill + border + left align),
# copy format from A1 which already carries that style.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Data row 2, text cells -> style index 2 (Text format "@" + left align)
$ws.Range("A2:E2").NumberFormat = "@"
$ws.Range("A2:E2").HorizontalAlignment = -4131
$ws.Range("G2:I2").NumberFormat = "@"
$ws.Range("G2:I2").HorizontalAlignment = -4131

# Data row 2, "code" style cells (numeric-ish text, no alignment) -> style index 4
$ws.Range("F2").NumberFormat = "@"
$ws.Range("J2:L2").NumberFormat = "@"

# Rows 3 & 4 placeholders
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").HorizontalAlignment = -4131
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").HorizontalAlignment = -4131

$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("B4").HorizontalAlignment = -4131
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").HorizontalAlignment = -4131
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").HorizontalAlignment = -4131
